# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail, same shape as the
# existing "2021-Q2" sheet) positioned right before the "总计" (summary)
# sheet, and prepends a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet immediately before "总计".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zongji)
$q1.Name = "2022-Q1"
# re-resolve after the structural change (Add shifts sheet indices)
$zongji = $wb.Worksheets.Item("总计")

# Match the look & feel (outline / page margins) of the other data sheets.
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Header row (style copied from the "总计" header so it matches s=2).
$zongji.Range("B1").Copy($q1.Range("B1:H1"))
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column A (style copied from "总计" A2 so it matches s=2).
$zongji.Range("A2").Copy($q1.Range("A2:A5"))
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2
$q1.Range("A5").Value = 3

# Columns B-G hold text (fund code must keep its leading zero, and the
# numbers are stored as plain text in the source workbook) - force text
# format before writing, then drop back to the default "Normal" style so
# no extra formatting is left behind on the cells.
$textRange = $q1.Range("B2:G5")
$textRange.NumberFormat = "@"

$q1.Range("B2").Value = "012027"
$q1.Range("C2").Value = "光大保德信安阳一年持有期混合型证券投资基金A"
$q1.Range("D2").Value = "15.22"
$q1.Range("E2").Value = "22.05"
$q1.Range("F2").Value = "0.80"
$q1.Range("G2").Value = "0.1218"

$q1.Range("B3").Value = "012010"
$q1.Range("C3").Value = "富国泰享回报6个月持有期混合型证券投资基金A"
$q1.Range("D3").Value = "9.29"
$q1.Range("E3").Value = "29.91"
$q1.Range("F3").Value = "1.00"
$q1.Range("G3").Value = "0.0929"

$q1.Range("B4").Value = "012028"
$q1.Range("C4").Value = "光大保德信安阳一年持有期混合型证券投资基金C"
$q1.Range("D4").Value = "7.68"
$q1.Range("E4").Value = "22.05"
$q1.Range("F4").Value = "0.80"
$q1.Range("G4").Value = "0.0614"

$q1.Range("B5").Value = "012011"
$q1.Range("C5").Value = "富国泰享回报6个月持有期混合型证券投资基金C"
$q1.Range("D5").Value = "0.09"
$q1.Range("E5").Value = "29.91"
$q1.Range("F5").Value = "1.00"
$q1.Range("G5").Value = "0.0009"

$textRange.Style = "Normal"

# Column H holds the numeric rank.
$q1.Range("H2").Value = 6
$q1.Range("H3").Value = 7
$q1.Range("H4").Value = 6
$q1.Range("H5").Value = 7

# ---------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计" (existing row 2 becomes
#    row 3, style untouched).
# ---------------------------------------------------------------------
$zongji.Rows(2).Insert()
# restore the index-column style on the newly inserted row (copied from
# the row pushed down to A3, which still carries the original s=2 style)
$zongji.Range("A3").Copy($zongji.Range("A2"))

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 4
$zongji.Range("D2").Value = 0.28
$zongji.Range("B2:D2").Style = "Normal"
# the 0-based row index in column A renumbers to account for the new row
$zongji.Range("A3").Value = 1

# ---------------------------------------------------------------------
# 3. Keep "2021-Q2" as the active/selected sheet, as it was originally.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q2").Activate()
